# Problema 2 probably working, need to check with teacher
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- b (B4) updated ---
$ws.Range("B4").Value = 10

# --- x table (C20:I26) tweaks ---
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 1

# --- t table (C28:I34) recalculated values ---
$ws.Range("D28").Value = 12
$ws.Range("E28").Value = 24
$ws.Range("F28").Value = 36
$ws.Range("G28").Value = 46
$ws.Range("H28").Value = 48
$ws.Range("I28").Value = 52

$ws.Range("D29").Value = 30
$ws.Range("E29").Value = 40
$ws.Range("F29").Value = 51
$ws.Range("G29").Value = 53
$ws.Range("H29").Value = 56
$ws.Range("I29").Value = 59

$ws.Range("C30").Value = 36
$ws.Range("D30").Value = 48
$ws.Range("E30").Value = 58
$ws.Range("F30").Value = 60
$ws.Range("G30").Value = 62
$ws.Range("H30").Value = 66
$ws.Range("I30").Value = 68

$ws.Range("C31").Value = 54
$ws.Range("D31").Value = 64
$ws.Range("E31").Value = 68
$ws.Range("F31").Value = 70
$ws.Range("G31").Value = 72
$ws.Range("H31").Value = 76
$ws.Range("I31").Value = 78

$ws.Range("C32").Value = 73
$ws.Range("D32").Value = 73
$ws.Range("E32").Value = 77
$ws.Range("F32").Value = 79
$ws.Range("G32").Value = 83
$ws.Range("H32").Value = 86
$ws.Range("I32").Value = 89

$ws.Range("C33").Value = 84
$ws.Range("D33").Value = 84
$ws.Range("E33").Value = 86
$ws.Range("F33").Value = 90
$ws.Range("G33").Value = 92
$ws.Range("H33").Value = 94
$ws.Range("I33").Value = 97

$ws.Range("C34").Value = 93
$ws.Range("D34").Value = 95
$ws.Range("E34").Value = 96
$ws.Range("F34").Value = 100
$ws.Range("G34").Value = 102
$ws.Range("H34").Value = 104
$ws.Range("I34").Value = 106

# --- selection moved to B5 ---
$ws.Range("B5").Select()
